$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the obsolete worker record row (row 18: 64915540 / MARTA LIGIA RUIZ ABAD / period 2508),
# shifting the rows below it up.
$ws.Rows("18:18").Delete()

# Center-align the "Periodo Mora" values for the remaining account-statement rows.
$ws.Range("E16:E17").HorizontalAlignment = -4108

# The deleted row used to close the table with a solid bottom border; re-apply that
# bottom border to the new last row of the table (now row 17).
$ws.Range("B17:J17").Borders.Item(9).LineStyle = 1
$ws.Range("B17:J17").Borders.Item(9).Weight = 2
$ws.Range("B17:J17").Borders.Item(9).ColorIndex = 1

# Update the totals to reflect the removal of that worker's debt record.
$ws.Range("E11").Value = 144000
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 2
